$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

# Row 51: San Telmo / Capital Sur
$ws.Range("M51").Value = -58.404058
$ws.Range("N51").Value = -34.634341
$ws.Range("O51").Value = "San Telmo"
$ws.Range("P51").Value = "Capital Sur"

# Row 52: Saavedra / Capital Norte
$ws.Range("M52").Value = -58.487821
$ws.Range("N52").Value = -34.554603
$ws.Range("O52").Value = "Saavedra"
$ws.Range("P52").Value = "Capital Norte"

# Row 53: No ubicado / No clasificado
$ws.Range("O53").Value = "No ubicado"
$ws.Range("P53").Value = "No clasificado, consultar con mantenimiento"

# Row 54: No ubicado / No clasificado
$ws.Range("O54").Value = "No ubicado"
$ws.Range("P54").Value = "No clasificado, consultar con mantenimiento"
